$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 333, shifting existing rows 333:348 down to 334:349
$ws.Rows(333).Insert()

# Populate the newly inserted row 333 with the new weekly data point
$ws.Range("A333").Value = 5
$ws.Range("B333").Value = 'Macroferia Regional de Talca'
$ws.Range("C333").Value = 'Maule'
$ws.Range("D333").Value = 45147
$ws.Range("E333").Value = 7
$ws.Range("F333").Value = 100112017
$ws.Range("G333").Value = 'Apio'
$ws.Range("H333").Value = 'Americana (o)'
$ws.Range("I333").Value = 'Primera'
$ws.Range("J333").Value = 800
$ws.Range("K333").Value = 5000
$ws.Range("L333").Value = 5000
$ws.Range("M333").Value = 5000
$ws.Range("N333").Value = '$/docena de matas'
$ws.Range("O333").Value = 'Provincia del Elquí'
$ws.Range("P333").Value = 833
$ws.Range("Q333").Value = 6
$ws.Range("R333").Value = 'Hortaliza'
